$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $find"
    }
}

Replace-Text "2024-08-03 Saturday" "2024-08-04 Sunday"
Replace-Text "461×8=3688" "790×7=5530"
Replace-Text "403×9=3627" "910×2=1820"
Replace-Text "669×9=6021" "192×9=1728"
Replace-Text "870×7=6090" "250×3=750"
Replace-Text "606×5=3030" "754×7=5278"
Replace-Text "976×6=5856" "216×8=1728"
Replace-Text "463×4=1852" "332×6=1992"
Replace-Text "826×9=7434" "404×5=2020"
Replace-Text "730×7=5110" "248×6=1488"
Replace-Text "322×3=966" "678×4=2712"
Replace-Text "961×9=8649" "851×3=2553"
Replace-Text "685×2=1370" "601×5=3005"
Replace-Text "759×7=5313" "713×5=3565"
Replace-Text "647×2=1294" "382×4=1528"
Replace-Text "311×6=1866" "956×9=8604"
Replace-Text "925×7=6475" "196×6=1176"
Replace-Text "545×4=2180" "702×2=1404"
Replace-Text "922×7=6454" "887×8=7096"
Replace-Text "617×9=5553" "524×2=1048"
Replace-Text "527×4=2108" "453×8=3624"
Replace-Text "212×7=1484" "316×3=948"
Replace-Text "434×6=2604" "616×6=3696"
Replace-Text "488×7=3416" "652×3=1956"
Replace-Text "257×2=514" "606×9=5454"
Replace-Text "946×7=6622" "654×6=3924"
